# "buy and sell flow added"
# Add a new "VendorEmail" worksheet at the end of the workbook, populate it
# with a header + a hyperlinked vendor e-mail address, and make it the
# active sheet (mirrors Excel's behaviour after inserting & selecting a
# new sheet interactively).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "VendorEmail"

# Row 1: label, Row 2: the vendor e-mail (with a mailto hyperlink).
$ws.Range("A1").Value = "Vendor Email "
$ws.Range("A2").Value = "testone@yopmail.com"
$ws.Range("A2").Hyperlinks.Add($ws.Range("A2"), "mailto:testone@yopmail.com")

# Match the hyperlink look used elsewhere in the workbook (e.g. signupdata!C2).
$ws.Range("A2").Style = $wb.Worksheets.Item("ValidLoginsheet").Range("A2").Style

# Column A width.
$ws.Columns.Item(1).ColumnWidth = 21.5

# Leave the selection on row 3, like the authored sheet (the newly
# inserted sheet is already the active / selected tab at this point).
[void]$ws.Rows.Item(3).Select()
